$d = $word.ActiveDocument

# Replace employee full name
$d.Content.Find.Execute("CRISTINA M. IGNO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FELICITAS M. SUMAGUI", 2)

# Replace position/title
$d.Content.Find.Execute("Administrative Aide III", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Casual Employee", 2)

# Replace office name (only the mail-merge field result, not the textbox occurrences)
$d.Content.Find.Execute("Human Resource Management Office", $true, $false, $false, $false, $false,
                         $true, 1, $false, "City Social Welfare Development Office", 1)

# Replace monthly salary values (appears twice)
$d.Content.Find.Execute("13,419.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11,814.00", 2)

# Replace total leave credits values (appears twice)
$d.Content.Find.Execute("325.931", $true, $false, $false, $false, $false,
                         $true, 1, $false, "120.916", 2)

# Replace total leave benefits value
$d.Content.Find.Execute("210,778.87", $true, $false, $false, $false, $false,
                         $true, 1, $false, "68,843.35", 2)
